$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G1").Value = 0.12816055735067355
$ws.Range("H1").Value = 0.18375715967599693
$ws.Range("I1").Value = 0.12492995561981513

$ws.Range("G2").Value = 0.090655737310379486
$ws.Range("H2").Value = 0.62510601173592995
$ws.Range("I2").Value = 0.41973996405031877

$ws.Range("G3").Value = 0.81367253694917996
$ws.Range("H3").Value = 0.72448108008498235
$ws.Range("I3").Value = 0.63899740950181672

$ws.Range("G4").Value = 0.0097880777824103111
$ws.Range("H4").Value = 0.03186017212454903
$ws.Range("I4").Value = 0.31874223744050778
